$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values changed (sample-size row)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON): B2, D2 and E2 were deleted outright; C2 got a new value
$ws.Range("B2").Value = $null
$ws.Range("C2").Value = -6.5711422825917634
$ws.Range("D2").Value = $null
$ws.Range("E2").Value = $null

# Row 3 (STR): B3:E3 recalculated/tweaked values
$ws.Range("B3").Value = -7.7900079309787529
$ws.Range("C3").Value = 4.0937642453683907
$ws.Range("D3").Value = -6.4305823250474115
$ws.Range("E3").Value = 25.643932419873284

# Selection narrowed from the whole data block to the edited B1:E3 block
$ws.Range("B1:E3").Select()
